$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Match the bold/border/centered header style used by the existing header row (A1:AC1)
# by copying the formatting from an existing header cell before setting values.
$ws.Range("A1").Copy() | Out-Null
$ws.Range("AD1:AF1").PasteSpecial(-4122) | Out-Null  # xlPasteFormats

# Add header cells for the new columns AD, AE, AF (Wins, Losses, Ties)
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Fill team record values for every data/footer row (2 through 52)
for ($row = 2; $row -le 52; $row++) {
    $ws.Cells.Item($row, 30).Value = 75  # AD
    $ws.Cells.Item($row, 31).Value = 87  # AE
    $ws.Cells.Item($row, 32).Value = 0   # AF
}
